$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header summary figures ---
# "VALOR MORA" total moves from 237250 to 9490 (only one worker/period remains)
$ws.Range("E11").Value2 = 9490
# "Cant. Trabajadores" drops from 3 to 1
$ws.Range("C13").Value2 = 1
# "Cant. Periodos" drops from 2 to 1
$ws.Range("F13").Value2 = 1

# --- Remove the two duplicated employee blocks (rows 17-19) ---
# Rows 16-20 originally held 3 employees across 5 rows (JOAN x2, ANDERSON x2, PEDRO x1).
# The update keeps only PEDRO JOSE QUITIAN GUZMAN's single record.
$ws.Rows("17:19").Delete()

# After that delete, PEDRO's record (previously row 20) is now row 17.
# Pull its values up into row 16 (which keeps its own original formatting/style),
# then drop the now-empty duplicate row.
$b17 = $ws.Range("B17").Value2
$c17 = $ws.Range("C17").Value2
$d17 = $ws.Range("D17").Value2
$e17 = $ws.Range("E17").Value2
$f17 = $ws.Range("F17").Value2
$g17 = $ws.Range("G17").Value2

$ws.Range("B16").Value2 = $b17
$ws.Range("C16").Value2 = $c17
$ws.Range("D16").Value2 = $d17
$ws.Range("E16").Value2 = $e17
$ws.Range("F16").Value2 = $f17
$ws.Range("G16").Value2 = $g17

$ws.Rows("17:17").Delete()

# Column D ("Nombre Trabajador") can now be narrower since the longest
# remaining name is shorter than before.
$ws.Columns("D:D").ColumnWidth = 27.8
